$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.082.53'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.81%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.564.68'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.17%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.81'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.97'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.65%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.577'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.62'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0814'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.47'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.37%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.959.09'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.22%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.108'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -4.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.561.32'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.47%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.07'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.845'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.51%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.112.80'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.85'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +4.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.59'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.25%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.56'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.99%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '253.19'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.78%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.06'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.93%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.44'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '39.96'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.81%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.23'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.68%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '154.16'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.32%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.42'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +4.02%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.34%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.42%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.83%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.02'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.15%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.54%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +5.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.119'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '22.61'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.53%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.92'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.01%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.67%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.003.26'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.78%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.03'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '83.39'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.811.80'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '74.39'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.18%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.44%  '
